# Applies commit "Add A3, A5, A7-A12" to the Rysiman evidence workbook.
$wb = $excel.ActiveWorkbook

# --- A2: append two more transfer rows (row 4 and row 5) ---
$ws = $wb.Worksheets.Item("A2")
$ws.Range("A4").Value = "14068E1916518AC76C9FB4F6A5F99FE9810C527DECC142B19924D88954ED17E5"
$ws.Range("B4").Value = "rysimangon11"
$ws.Range("C4").Value = "rysigon44"
$ws.Range("A5").Value = "E4CCD456D8E1124C5E4B063A9983731F17ECAA509B7648D824357D49999CE496"
$ws.Range("B5").Value = "rysimangon11"
$ws.Range("C5").Value = "rysigon04"
$ws.Range("A4").Select()
$ws.Range("A4:C5").Select()

# --- A3: fill in the previously-empty evidence row ---
$ws = $wb.Worksheets.Item("A3")
$ws.Range("A2").Value = "03A69A159F1DAC5F0490027D0F096EB6E8431A60498D532993144440ED327AB8"
$ws.Range("B2").Value = "stars1mgcpkhw4yx4hhygtzt99wlj2d8el23g29g8x0f2zf754a0clcr0q4exy9m"
$ws.Range("C2").Value = "rysigon04"
$ws.Range("D2").Value = "elgafar-1"
$ws.Range("B2").Select()

# --- A5: replace the placeholder evidence row with real values ---
$ws = $wb.Worksheets.Item("A5")
$ws.Range("A2").Value = "75D8ECB14016DD5D89CE8EA56ACDB311228C5F8C848B36C808B64EC8F28B6843"
$ws.Range("B2").Value = "stars1mgcpkhw4yx4hhygtzt99wlj2d8el23g29g8x0f2zf754a0clcr0q4exy9m"
$ws.Range("C2").Value = "rysigon04"
$ws.Range("D2").Value = "elgafar-1"
$ws.Range("B2").Select()

# --- A7..A12: fill in ClassID / NFTID placeholders ---
$ws = $wb.Worksheets.Item("A7")
$ws.Range("A2").Value = "ibc/33C71002A5D3A4E329AF1004F3AC67214BCA5B7BE4C5B257C2A9FF2B3F1D7A6C"
$ws.Range("B2").Value = "rysigon07"
$ws.Range("C9").Select()

$ws = $wb.Worksheets.Item("A8")
$ws.Range("A2").Value = "ibc/33C71002A5D3A4E329AF1004F3AC67214BCA5B7BE4C5B257C2A9FF2B3F1D7A6C"
$ws.Range("B2").Value = "rysigon08"
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("A9")
$ws.Range("A2").Value = "ibc/33C71002A5D3A4E329AF1004F3AC67214BCA5B7BE4C5B257C2A9FF2B3F1D7A6C"
$ws.Range("B2").Value = "rysigon09"
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("A10")
$ws.Range("A2").Value = "ibc/33C71002A5D3A4E329AF1004F3AC67214BCA5B7BE4C5B257C2A9FF2B3F1D7A6C"
$ws.Range("B2").Value = "rysigon10"
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("A11")
$ws.Range("A2").Value = "ibc/33C71002A5D3A4E329AF1004F3AC67214BCA5B7BE4C5B257C2A9FF2B3F1D7A6C"
$ws.Range("B2").Value = "rysigon11"
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("A12")
$ws.Range("A2").Value = "ibc/33C71002A5D3A4E329AF1004F3AC67214BCA5B7BE4C5B257C2A9FF2B3F1D7A6C"
$ws.Range("B2").Value = "rysigon12"
$ws.Activate()
$ws.Range("I8").Select()
